$d = $word.ActiveDocument

$replacements = @(
    @("454÷6=75, 4", "748÷5=149, 3"),
    @("396÷7=56, 4", "315÷4=78, 3"),
    @("777÷3=259, 0", "339÷8=42, 3"),
    @("923÷2=461, 1", "812÷9=90, 2"),
    @("682÷6=113, 4", "356÷4=89, 0"),
    @("935÷9=103, 8", "344÷3=114, 2"),
    @("277÷5=55, 2", "862÷9=95, 7"),
    @("336÷7=48, 0", "263÷2=131, 1"),
    @("297÷7=42, 3", "857÷8=107, 1"),
    @("726÷9=80, 6", "410÷9=45, 5"),
    @("555÷8=69, 3", "744÷2=372, 0"),
    @("984÷7=140, 4", "588÷8=73, 4"),
    @("833÷3=277, 2", "949÷5=189, 4"),
    @("504÷6=84, 0", "945÷4=236, 1"),
    @("935÷8=116, 7", "695÷7=99, 2"),
    @("488÷3=162, 2", "400÷3=133, 1"),
    @("160÷5=32, 0", "721÷9=80, 1"),
    @("579÷6=96, 3", "184÷5=36, 4"),
    @("161÷9=17, 8", "464÷2=232, 0"),
    @("417÷9=46, 3", "430÷4=107, 2"),
    @("866÷9=96, 2", "897÷6=149, 3"),
    @("650÷4=162, 2", "177÷9=19, 6"),
    @("718÷4=179, 2", "597÷7=85, 2"),
    @("851÷4=212, 3", "852÷4=213, 0"),
    @("807÷7=115, 2", "864÷5=172, 4")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new, 2)
}

$d.Save()
